# Fix the typo "Enigneer" -> "Engineer" (and re-flow it as "DevOps Engineer")
# in the bio paragraph on slide 1, as described by the commit diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the bio textbox ("Rectangle 1") by name so this is resilient to any
# shape re-ordering.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 1") {
        $target = $candidate
    }
}

$tr = $target.TextFrame.TextRange

# The bio text lives in the third paragraph ("Introduction: ...").
$para = $tr.Paragraphs(3, 1)

# 1) " LLC as DevOps " -> " LLC as " (the word "DevOps" moves down to join
#    the corrected "Engineer" run).
$full = $para.Text
$i1 = $full.IndexOf(" LLC as DevOps ")
if ($i1 -ge 0) {
    $rng1 = $para.Characters($i1 + 1, (" LLC as DevOps ").Length)
    $rng1.Text = " LLC as "
}

# 2) "Enigneer" -> "DevOps Engineer" (fixes the misspelling).
$full2 = $para.Text
$i2 = $full2.IndexOf("Enigneer")
if ($i2 -ge 0) {
    $rng2 = $para.Characters($i2 + 1, ("Enigneer").Length)
    $rng2.Text = "DevOps Engineer"
}

# 3) Re-assert the remainder of the sentence so it settles back into a
#    single contiguous run.
$full3 = $para.Text
$tail = ", where I gained hands-on experience in DevOps and software engineering. I am passionate about leveraging technology to drive efficiency and innovation."
$i3 = $full3.IndexOf(", where I gained")
if ($i3 -ge 0) {
    $rng3 = $para.Characters($i3 + 1, $full3.Length - $i3)
    $rng3.Text = $tail
}
